# Append: 2026-01-26 01:45 JST
# Updates the "ランサーズ" (Lancers) sheet in 案件情報.xlsx:
#  - refresh the "取得日時" (fetched-at) timestamp on existing rows 2-10
#  - widen column D slightly
#  - row 7 (a new scrape of the same listing slot) gets new title/price/url/score/skill text
#  - append a brand-new row 11 with its own hyperlink

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2026-01-26 01:45:08"

# --- refresh "取得日時" (column A) for every existing data row ---
$ws.Range("A2").Value = $timestamp
$ws.Range("A3").Value = $timestamp
$ws.Range("A4").Value = $timestamp
$ws.Range("A5").Value = $timestamp
$ws.Range("A6").Value = $timestamp
$ws.Range("A7").Value = $timestamp
$ws.Range("A8").Value = $timestamp
$ws.Range("A9").Value = $timestamp
$ws.Range("A10").Value = $timestamp

# --- widen column D (28 -> 32 characters) ---
$ws.Range("D1").ColumnWidth = 31.17

# --- row 7 got replaced by a different listing on this re-scrape ---
$ws.Range("B7").Value = "自動化システム"
$ws.Range("D7").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5477084"
$ws.Range("G7").Value = 110
$ws.Range("H7").Value = "◆自動化"

# --- append new row 11 ---
$ws.Range("A11").Value = $timestamp
$ws.Range("B11").Value = "【医療保険】オンライン資格確認・請求端末セットアップ依頼"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5478715"
$ws.Range("G11").Value = 13

$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5478715")
$ws.Range("F11").Style = "Hyperlink"
